# Update gh-pages to output generated at 456a3b4
# Apply data refresh to the 合肥-漫展信息 workbook: update "想去人数" (F column)
# counts, and mark row 3's "最低票价" (G column) as sold out ("不可售") on the
# sheets that list it ("展览" and "全部类型"); bump the "演出" sheet's single
# row too.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 125
$ws1.Range("F3").Value = 233
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F5").Value = 6652
$ws1.Range("F8").Value = 136
$ws1.Range("F9").Value = 6106
$ws1.Range("F10").Value = 43
$ws1.Range("F15").Value = 92
$ws1.Range("F16").Value = 388
$ws1.Range("F19").Value = 358
$ws1.Range("F22").Value = 4431
$ws1.Range("F23").Value = 50
$ws1.Range("F24").Value = 23
$ws1.Range("F25").Value = 188
$ws1.Range("F26").Value = 41

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 44

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 125
$ws4.Range("F3").Value = 233
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F5").Value = 6652
$ws4.Range("F8").Value = 136
$ws4.Range("F9").Value = 6106
$ws4.Range("F10").Value = 43
$ws4.Range("F15").Value = 92
$ws4.Range("F16").Value = 388
$ws4.Range("F19").Value = 358
$ws4.Range("F22").Value = 4431
$ws4.Range("F23").Value = 44
$ws4.Range("F24").Value = 50
$ws4.Range("F25").Value = 23
$ws4.Range("F26").Value = 188
$ws4.Range("F27").Value = 41
